$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fun fact text in V1: "realname" -> "truename"
$ws.Range("V1").Value = "Pangil's truename is only known by a handful of animals."

# Update the selected cell to V1 (as shown in the diff's sheetView selection)
$ws.Range("V1").Select()
